$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '71.381.11'
$ws.Cells.Item(2, 5).Value = '  -2.41%  '
$ws.Cells.Item(3, 4).Value = '3.879.41'
$ws.Cells.Item(3, 5).Value = '  -2.86%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.998'
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '608.52'
$ws.Cells.Item(5, 5).Value = '  +0.75%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '171.94'
$ws.Cells.Item(6, 5).Value = '  +4.98%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.672'
$ws.Cells.Item(7, 5).Value = '  -2.07%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  +0.03%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.752'
$ws.Cells.Item(9, 5).Value = '  -0.37%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.179'
$ws.Cells.Item(10, 5).Value = '  +5.78%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '53.91'
$ws.Cells.Item(11, 5).Value = '  -1.57%  '
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000324'
$ws.Cells.Item(12, 5).Value = '  +1.03%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.57'
$ws.Cells.Item(13, 5).Value = '  +4.69%  '
$ws.Cells.Item(14, 4).Value = '4.484.47'
$ws.Cells.Item(14, 5).Value = '  -3.10%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '21.25'
$ws.Cells.Item(15, 5).Value = '  +3.30%  '
$ws.Cells.Item(16, 4).Value = '3.870.75'
$ws.Cells.Item(16, 5).Value = '  -3.17%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.00'
$ws.Cells.Item(17, 5).Value = '  -1.15%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '1.22'
$ws.Cells.Item(18, 5).Value = '  -3.90%  '
$ws.Cells.Item(19, 5).Value = '  -2.20%  '
$ws.Cells.Item(20, 4).Value = '71.091.56'
$ws.Cells.Item(20, 5).Value = '  -2.35%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '441.75'
$ws.Cells.Item(21, 5).Value = '  +0.22%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.83'
$ws.Cells.Item(22, 5).Value = '  +0.63%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '95.00'
$ws.Cells.Item(23, 5).Value = '  -2.03%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '3.31'
$ws.Cells.Item(24, 5).Value = '  -4.38%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '13.92'
$ws.Cells.Item(25, 5).Value = '  -3.43%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '11.79'
$ws.Cells.Item(26, 5).Value = '  +3.38%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '4.03'
$ws.Cells.Item(27, 5).Value = '  -7.09%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '5.97'
$ws.Cells.Item(28, 5).Value = '  +0.14%  '
$ws.Cells.Item(29, 5).Value = '  +0.92%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '8.78'
$ws.Cells.Item(30, 5).Value = '  +10.97%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '35.29'
$ws.Cells.Item(31, 5).Value = '  -3.33%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '13.60'
$ws.Cells.Item(32, 5).Value = '  -2.33%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '48.04'
$ws.Cells.Item(33, 5).Value = '  -2.23%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.127'
$ws.Cells.Item(34, 5).Value = '  -3.89%  '
$ws.Cells.Item(35, 4).Value = '0.0₃0996'
$ws.Cells.Item(35, 5).Value = '  +10.65%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '69.19'
$ws.Cells.Item(36, 5).Value = '  -2.34%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '639.64'
$ws.Cells.Item(37, 5).Value = '  -3.80%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.443'
$ws.Cells.Item(38, 5).Value = '  +1.01%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.148'
$ws.Cells.Item(39, 5).Value = '  +0.58%  '
$ws.Cells.Item(40, 5).Value = '  +0.15%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.00'
$ws.Cells.Item(41, 5).Value = '  -0.05%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.27'
$ws.Cells.Item(42, 5).Value = '  -3.15%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.89'
$ws.Cells.Item(43, 5).Value = '  +9.10%  '
$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '3.18'
$ws.Cells.Item(44, 5).Value = '  +18.30%  '
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0473'
$ws.Cells.Item(45, 5).Value = '  -3.58%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '10.23'
$ws.Cells.Item(46, 5).Value = '  -4.72%  '
$ws.Cells.Item(47, 2).Value = 'WEMIXToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.92'
$ws.Cells.Item(47, 5).Value = '  -12.91%  '
$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.144'
$ws.Cells.Item(48, 5).Value = '  -4.00%  '
$ws.Cells.Item(49, 4).Value = '2.942.26'
$ws.Cells.Item(49, 5).Value = '  +1.01%  '
$ws.Cells.Item(50, 5).Value = '  -3.40%  '
$ws.Cells.Item(51, 5).Value = '  +3.20%  '
